$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: In the "Functional Requirements Analysis Tables" table, the cell
# that holds "R1:" gets a new run appended: " Customize character" (Arial).
# ---------------------------------------------------------------------------
$reqTable = $d.Tables.Item(2)
$r1Cell = $reqTable.Cell(1, 2)
$r1Para = $r1Cell.Range.Paragraphs.Item(1)
$r1ParaRange = $r1Para.Range

# Target only the run's text (exclude the trailing paragraph mark) so the
# insertion replaces just the "R1:" run with "R1:" + the new run, keeping
# the paragraph's own properties (pPr) untouched.
$r1TextRange = $d.Range($r1ParaRange.Start, $r1ParaRange.End - 1)

$r1Xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="37BB83E2" w14:textId="1E3E0632" w:rsidR="3DB678F9" w:rsidRDefault="3DB678F9" w:rsidP="3DB678F9"><w:pPr><w:widowControl w:val="0"/><w:spacing w:before="240" w:line="276" w:lineRule="auto"/><w:ind w:left="144"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr><w:r w:rsidRPr="3DB678F9"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>R1:</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> Customize character</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r1TextRange.InsertXML($r1Xml)

# ---------------------------------------------------------------------------
# Edit 2: In the same table, the "Summary" row's value cell holds a single
# empty paragraph whose pPr spacing/indent changes:
#   <w:spacing w:before="120"/>  ->  <w:spacing w:before="120" w:after="120"/>
#   (new)                        ->  <w:ind w:left="144" w:right="144"/>
# ---------------------------------------------------------------------------
$summaryCell = $reqTable.Cell(2, 2)
$summaryPara = $summaryCell.Range.Paragraphs.Item(1)
$summaryParaRange = $summaryPara.Range

$summaryXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="3EB8B502" w14:textId="5A6A0A85" w:rsidR="3DB678F9" w:rsidRDefault="3DB678F9" w:rsidP="3DB678F9"><w:pPr><w:spacing w:before="120" w:after="120"/><w:ind w:left="144" w:right="144"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$summaryParaRange.InsertXML($summaryXml)
